# Automatische test-sync: 2025-08-03 18:29:50
# Appends a new logged test-mail row to the "Logs" sheet and refreshes the
# corresponding count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

# Next free row right after the existing data (row 35 -> row 36).
$newRow = $logs.Cells.Item($logs.Rows.Count, 1).End(-4162).Row + 1

$logs.Cells.Item($newRow, 1).Value = "Kun je nagaan of dit nog leverbaar is?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #8: Kun je nagaan of dit nog leverbaar is?"
$logs.Cells.Item($newRow, 4).Value = "Inkoop / Bestellingen"
$logs.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value = "2025-08-03 18:28:51"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# Update the Dashboard summary count for "Inkoop / Bestellingen" (row 4, col B).
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(4, 2).Value = 8

# Extend the conditional-formatting coverage of each column so it keeps
# including the whole data range through the newly appended row.
$logs.Range("D2:D35").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D36"))
$logs.Range("G2:G35").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G36"))
$logs.Range("H2:H35").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H36"))
$logs.Range("I2:I35").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I36"))
$logs.Range("J2:J35").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J36"))
